# Apply scheduled-runner profit recalculations to the Titan_Profits sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("N3").ClearContents()
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("K98").Value = 450907.03
$ws.Range("N98").Value = -11212
$ws.Range("H98").Value = 433880.47
$ws.Range("M98").Value = -449409.03
$ws.Range("L98").Value = 8216
$ws.Range("J98").Value = 8216
$ws.Range("I98").Value = 450907.03
$ws.Range("H102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("N102").ClearContents()
$ws.Range("L122").Value = 24648
$ws.Range("N122").Value = -29548
$ws.Range("M122").Value = -1350271.09
$ws.Range("K122").Value = 1352721.09
$ws.Range("H122").Value = 433880.47
$ws.Range("I122").Value = 450907.03
$ws.Range("J122").Value = 8216

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 19416
$ws.Range("I32").Value = 4625.0586
$ws.Range("K32").Value = 4625.0586
$ws.Range("M32").Value = -4338.0586
$ws.Range("L122").Value = 7821
$ws.Range("N122").Value = -12721
$ws.Range("M122").Value = -23481.118
$ws.Range("K122").Value = 25931.118
$ws.Range("H122").Value = 8008.263
$ws.Range("I122").Value = 8643.706
$ws.Range("J122").Value = 2607
$ws.Range("J130").Value = 25214.5
$ws.Range("H130").Value = 25214.5
$ws.Range("L130").Value = 25214.5
$ws.Range("N130").Value = -35254.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("K86").Value = 1340.625
$ws.Range("L86").Value = 1437.375
$ws.Range("I86").Value = 1340.625
$ws.Range("N86").Value = -3683.375
$ws.Range("J86").Value = 1437.375
$ws.Range("H86").Value = 1389
$ws.Range("M86").Value = -217.625
$ws.Range("J89").Value = 1437.375
$ws.Range("K89").Value = 6703.125
$ws.Range("H89").Value = 1389
$ws.Range("I89").Value = 1340.625
$ws.Range("L89").Value = 7186.875
$ws.Range("M89").Value = -1087.125
$ws.Range("N89").Value = -18418.875
$ws.Range("H94").Value = 831.375
$ws.Range("K94").Value = 714.2143
$ws.Range("I94").Value = 714.2143
$ws.Range("L94").Value = 995.4
$ws.Range("J94").Value = 995.4
$ws.Range("M94").Value = -263.2143
$ws.Range("N94").Value = -1897.4
$ws.Range("J99").Value = 2666.3333
$ws.Range("L99").Value = 2666.3333
$ws.Range("H99").Value = 1880.5385
$ws.Range("N99").Value = -5662.3333
$ws.Range("K99").Value = 1644.8
$ws.Range("I99").Value = 1644.8
$ws.Range("M99").Value = -146.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("J31").Value = 7242.6665
$ws.Range("N31").Value = -7832.6665
$ws.Range("M31").Value = -1072.5938
$ws.Range("I31").Value = 1367.5938
$ws.Range("L31").Value = 7242.6665
$ws.Range("K31").Value = 1367.5938
$ws.Range("H31").Value = 3242.617
$ws.Range("L34").Value = 7242.6665
$ws.Range("J34").Value = 7242.6665
$ws.Range("N34").Value = -7646.6665
$ws.Range("K34").Value = 1367.5938
$ws.Range("H34").Value = 3242.617
$ws.Range("I34").Value = 1367.5938
$ws.Range("M34").Value = -1165.5938
$ws.Range("I62").Value = 19680.285
$ws.Range("N62").Value = -8738
$ws.Range("K62").Value = 19680.285
$ws.Range("J62").Value = 7490
$ws.Range("L62").Value = 7490
$ws.Range("M62").Value = -19056.285
$ws.Range("H62").Value = 16472.316
$ws.Range("K65").Value = 98401.425
$ws.Range("H65").Value = 16472.316
$ws.Range("N65").Value = -43690
$ws.Range("I65").Value = 19680.285
$ws.Range("J65").Value = 7490
$ws.Range("L65").Value = 37450
$ws.Range("M65").Value = -95281.425
$ws.Range("J99").Value = 3155.5557
$ws.Range("L99").Value = 3155.5557
$ws.Range("H99").Value = 2629
$ws.Range("N99").Value = -6151.5557
$ws.Range("K99").Value = 2102.4443
$ws.Range("I99").Value = 2102.4443
$ws.Range("M99").Value = -604.4443000000001
$ws.Range("K105").Value = 1126
$ws.Range("M105").Value = 621
$ws.Range("I105").Value = 1126
$ws.Range("H105").Value = 1105
$ws.Range("N105").Value = -4494
$ws.Range("L105").Value = 1000
$ws.Range("J105").Value = 1000
$ws.Range("H107").Value = 754
$ws.Range("M107").Value = 1472.5
$ws.Range("L107").Value = 1980
$ws.Range("I107").Value = 447.5
$ws.Range("J107").Value = 1980
$ws.Range("N107").Value = -5820
$ws.Range("K107").Value = 447.5
$ws.Range("N126").Value = -14406.6671
$ws.Range("K126").Value = 6307.3329
$ws.Range("L126").Value = 9466.667099999999
$ws.Range("H126").Value = 2629
$ws.Range("J126").Value = 3155.5557
$ws.Range("I126").Value = 2102.4443
$ws.Range("M126").Value = -3837.3329

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("N3").Value = -18224
$ws.Range("K3").Value = 6000
$ws.Range("M3").Value = -5888
$ws.Range("J3").Value = 6000
$ws.Range("L3").Value = 18000
$ws.Range("H3").Value = 5555.5557
$ws.Range("I3").Value = 2000
$ws.Range("H92").Value = 2260
$ws.Range("J92").Value = 2260
$ws.Range("N92").Value = -9276
$ws.Range("L92").Value = 6780
$ws.Range("H101").Value = 8000
$ws.Range("J101").Value = 0
$ws.Range("M101").Value = -21566
$ws.Range("I101").Value = 8000
$ws.Range("K101").Value = 24000
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()
$ws.Range("L122").Value = 20933.1
$ws.Range("N122").Value = -25833.1
$ws.Range("M122").Value = -1150
$ws.Range("K122").Value = 3600
$ws.Range("H122").Value = 1532.8823
$ws.Range("I122").Value = 400
$ws.Range("J122").Value = 2325.9

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("L5").Value = 1000
$ws.Range("J5").Value = 1000
$ws.Range("N5").Value = -1224
$ws.Range("H5").Value = 1000
$ws.Range("H24").Value = 19642.857
$ws.Range("L24").Value = 7500
$ws.Range("J24").Value = 7500
$ws.Range("N24").Value = -7846
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()
$ws.Range("M122").Value = -4882
$ws.Range("K122").Value = 7332
$ws.Range("H122").Value = 2444
$ws.Range("I122").Value = 2444
$ws.Range("J122").Value = 0
$ws.Range("N126").Value = -13818.8
$ws.Range("K126").Value = 7322.625
$ws.Range("L126").Value = 8878.799999999999
$ws.Range("H126").Value = 2779.1738
$ws.Range("J126").Value = 2959.6
$ws.Range("I126").Value = 2440.875
$ws.Range("M126").Value = -4852.625

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3213.2354
$ws.Range("I7").Value = 2644
$ws.Range("M7").Value = -2532
$ws.Range("J7").Value = 3450.4167
$ws.Range("L7").Value = 3450.4167
$ws.Range("K7").Value = 2644
$ws.Range("N7").Value = -3674.4167
$ws.Range("N126").Value = -15291.2501
$ws.Range("K126").Value = 7932
$ws.Range("L126").Value = 10351.2501
$ws.Range("H126").Value = 3213.2354
$ws.Range("J126").Value = 3450.4167
$ws.Range("I126").Value = 2644
$ws.Range("M126").Value = -5462

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("N81").Value = -10966.444
$ws.Range("M81").Value = -2806.8
$ws.Range("K81").Value = 3867.8
$ws.Range("L81").Value = 8844.444
$ws.Range("H81").Value = 3112.5789
$ws.Range("I81").Value = 1933.9
$ws.Range("J81").Value = 4422.222
$ws.Range("N84").Value = -54830.22
$ws.Range("J84").Value = 4422.222
$ws.Range("H84").Value = 3112.5789
$ws.Range("I84").Value = 1933.9
$ws.Range("L84").Value = 44222.22
$ws.Range("M84").Value = -14035
$ws.Range("K84").Value = 19339
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()
